$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-5) are reordered by date. Update the D (Fecha),
# M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio
# ponderado) and S (Precio $/Kg) columns to reflect the new row order.

# Row 2 -> old Row 5 values
$ws.Range("D2").Value = 44257
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 806

# Row 3 -> old Row 2 values
$ws.Range("D3").Value = 44252
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 13500
$ws.Range("S3").Value = 750

# Row 4 -> old Row 3 values
$ws.Range("D4").Value = 44250
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 806

# Row 5 -> old Row 4 values
$ws.Range("D5").Value = 44253
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 806
